$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: " All new RDS instances since launch use this integration by
# default when rotation is configured on any database users." becomes a
# multi-run sentence with an inserted "(~50k/year)" aside and "on DB users"
# instead of "on any database users".
# -----------------------------------------------------------------------
$oldSentence = " All new RDS instances since launch use this integration by default when rotation is configured on any database users."
$find1 = $d.Content
$find1.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target1 = $d.Range($find1.Start, $find1.End)

$newXml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:eastAsia="Times New Roman" w:hAnsi="Garamond" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t xml:space="preserve"> All new RDS instances since launch </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:eastAsia="Times New Roman" w:hAnsi="Garamond" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t>(</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:eastAsia="Times New Roman" w:hAnsi="Garamond" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t>~</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:eastAsia="Times New Roman" w:hAnsi="Garamond" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t xml:space="preserve">50k/year) </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:eastAsia="Times New Roman" w:hAnsi="Garamond" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t>use this integration by default when rotation is configured</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:eastAsia="Times New Roman" w:hAnsi="Garamond" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t xml:space="preserve"> on DB users</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Garamond" w:eastAsia="Times New Roman" w:hAnsi="Garamond" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t>.</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target1.InsertXML($newXml1)

# -----------------------------------------------------------------------
# Change 2: inside the hyperlink "Redshift's Service-Linked Admin Secrets
# feature", "Service-Linked" becomes "Auto-Created" and the run fragments
# are collapsed into a single run (text becomes
# "Redshift's Auto-Created Admin Secrets feature").
# -----------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("Service-Linked", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target2 = $d.Range($find2.Start, $find2.End)
$target2.Text = "Auto-Created"
